# Automatic update of files.
# Rotates the "Tretåig hackspett" / "Garnlav" observation records across
# rows 11-15 and swaps rows 19-20, per the upstream artfynd export refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11 ----
$ws.Range("A11").Value = 131167653
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "1"
$ws.Range("L11").Value = "hona"
$ws.Range("M11").Value = "födosökande"
$ws.Range("N11").Value = "observerad"
$ws.Range("Q11").Value = 613336
$ws.Range("R11").Value = 6997445
$ws.Range("Z11").Value = "13:47"
$ws.Range("AB11").Value = "13:48"
$ws.Range("AJ11").Value = ""
$ws.Range("AK11").Value = ""
$ws.Range("AO11").Value = ""

# ---- Row 12 ----
$ws.Range("A12").Value = 131167651
$ws.Range("I12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "äldre spår"
$ws.Range("N12").Value = ""
$ws.Range("Q12").Value = 613285
$ws.Range("R12").Value = 6997537
$ws.Range("Z12").Value = ""
$ws.Range("AB12").Value = ""
$ws.Range("AC12").Value = "Äldre ringhack på tall"

# ---- Row 13 ----
$ws.Range("A13").Value = 131167669
$ws.Range("B13").Value = 79244
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("M13").Value = ""
$ws.Range("Q13").Value = 613256
$ws.Range("R13").Value = 6997380
$ws.Range("AC13").Value = ""
$ws.Range("AJ13").Value = "tall"
$ws.Range("AK13").Value = "Pinus sylvestris"
$ws.Range("AO13").Value = "Pinus sylvestris"

# ---- Row 14 ----
$ws.Range("A14").Value = 131167650
$ws.Range("B14").Value = 57884
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = "Tretåig hackspett"
$ws.Range("G14").Value = "Picoides tridactylus"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("M14").Value = "färska spår"
$ws.Range("Q14").Value = 613278
$ws.Range("R14").Value = 6997506
$ws.Range("AC14").Value = "Färska ringhack  på tall"

# ---- Row 15 ----
$ws.Range("A15").Value = 131167670
$ws.Range("B15").Value = 79244
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("M15").Value = ""
$ws.Range("Q15").Value = 613271
$ws.Range("R15").Value = 6997395
$ws.Range("AC15").Value = "Observerad på tre granar inom 10m radie"

# ---- Row 19 ----
$ws.Range("A19").Value = 131167652
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = "1"
$ws.Range("M19").Value = "födosökande"
$ws.Range("N19").Value = "observerad"
$ws.Range("Q19").Value = 613264
$ws.Range("R19").Value = 6997532
$ws.Range("Z19").Value = "12:46"
$ws.Range("AB19").Value = "12:48"
$ws.Range("AC19").Value = "Hackspetten syns i profil på ett smalt träd i centrum av den tagna bilden."

# ---- Row 20 ----
$ws.Range("A20").Value = 131167654
$ws.Range("I20").Value = ""
$ws.Range("M20").Value = "äldre spår"
$ws.Range("N20").Value = ""
$ws.Range("Q20").Value = 613254
$ws.Range("R20").Value = 6997565
$ws.Range("Z20").Value = ""
$ws.Range("AB20").Value = ""
$ws.Range("AC20").Value = "Äldre ringhack på tall"
